$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row (C2:C348).
# All of them currently equal 45177 and must be bumped to 45178 (i.e. +1 day).
$ws.Range("C2:C348").Value = 45178
